# Apply the "single child" edit to the g__Kandleria worksheet:
#  - Drop the "max" column (old column C) entirely.
#  - Rename header cell (now C1) from "max"/"prediction" shuffle to "prediction".
#  - Keep "rejection-f" header, now in column D.
#  - Collapse the two data rows (RUG287.fasta, RUG655.fasta) down to a
#    single data row for RUG287.fasta, with an updated numeric value in B2.
#  - Shrink the used range down to A1:D2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old column C ("max") entirely, shifting D/E left ---
$ws.Range("C:C").Delete()

# --- Row 1 headers (now: A=Row, B=1-s__..., C=prediction, D=rejection-f) ---
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "1-s__Kandleria vitulina"
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# --- Row 2 data ---
$ws.Range("A2").Value = "RUG287.fasta"
$ws.Range("B2").Value = 11.97264756875075
$ws.Range("C2").Value = "s__Kandleria vitulina"
$ws.Range("D2").Value = "s__Kandleria vitulina"

# --- Remove row 3 (the RUG655.fasta row) entirely ---
$ws.Range("A3:D3").Delete()

# Tidy up: drop anything lingering beyond the new D2 extent.
$ws.Range("E1:E3").Clear()
$ws.Range("A4:D4").Clear()
